$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp (A1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 15 de Julio de 2020 a las 12:52"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Cells.Item(4,2).Value = 3546278
$ws.Cells.Item(4,3).Value = 1201
$ws.Cells.Item(4,5).Value = 1806206
$ws.Cells.Item(4,7).Value = 19
$ws.Cells.Item(4,8).Value = 139162

# Row 6: 'India' -> 'India'
$ws.Cells.Item(6,2).Value = 939192
$ws.Cells.Item(6,3).Value = 1705
$ws.Cells.Item(6,4).Value = 593198
$ws.Cells.Item(6,5).Value = 321667

# Row 14: 'Iran' -> 'Iran'
$ws.Cells.Item(14,2).Value = 264561
$ws.Cells.Item(14,3).Value = 2388
$ws.Cells.Item(14,4).Value = 227561
$ws.Cells.Item(14,5).Value = 23590
$ws.Cells.Item(14,7).Value = 199
$ws.Cells.Item(14,8).Value = 13410

# Row 34: 'Belgica' -> 'Belgica'
$ws.Cells.Item(34,2).Value = 62872
$ws.Cells.Item(34,3).Value = 91
$ws.Cells.Item(34,4).Value = 17242
$ws.Cells.Item(34,5).Value = 35842
$ws.Cells.Item(34,7).Value = 1
$ws.Cells.Item(34,8).Value = 9788

# Row 35: 'Oman' -> 'Oman'
$ws.Cells.Item(35,2).Value = 61247
$ws.Cells.Item(35,3).Value = 1679
$ws.Cells.Item(35,4).Value = 39038
$ws.Cells.Item(35,5).Value = 21928
$ws.Cells.Item(35,7).Value = 8
$ws.Cells.Item(35,8).Value = 281

# Row 36: 'Filipinas' -> 'Filipinas'
$ws.Cells.Item(36,2).Value = 58850
$ws.Cells.Item(36,3).Value = 1392
$ws.Cells.Item(36,4).Value = 20976
$ws.Cells.Item(36,5).Value = 36260
$ws.Cells.Item(36,7).Value = 11
$ws.Cells.Item(36,8).Value = 1614

# Row 49: 'Barein' -> 'Rumania'
$ws.Cells.Item(49,1).Value = "Rumania"
$ws.Cells.Item(49,2).Value = 34226
$ws.Cells.Item(49,3).Value = 641
$ws.Cells.Item(49,4).Value = 22049
$ws.Cells.Item(49,5).Value = 10225
$ws.Cells.Item(49,7).Value = 21
$ws.Cells.Item(49,8).Value = 1952

# Row 50: 'Nigeria' -> 'Barein'
$ws.Cells.Item(50,1).Value = "Barein"
$ws.Cells.Item(50,2).Value = 34078
$ws.Cells.Item(50,4).Value = 29753
$ws.Cells.Item(50,5).Value = 4211
$ws.Cells.Item(50,7).Value = 3
$ws.Cells.Item(50,8).Value = 114

# Row 51: 'Rumania' -> 'Nigeria'
$ws.Cells.Item(51,1).Value = "Nigeria"
$ws.Cells.Item(51,2).Value = 33616
$ws.Cells.Item(51,4).Value = 13792
$ws.Cells.Item(51,5).Value = 19070
$ws.Cells.Item(51,8).Value = 754

# Row 57: 'Azerbaiyan' -> 'Ghana'
$ws.Cells.Item(57,1).Value = "Ghana"
$ws.Cells.Item(57,2).Value = 25252
$ws.Cells.Item(57,3).Value = 264
$ws.Cells.Item(57,4).Value = 21397
$ws.Cells.Item(57,5).Value = 3716
$ws.Cells.Item(57,8).Value = 139

# Row 58: 'Ghana' -> 'Azerbaiyan'
$ws.Cells.Item(58,1).Value = "Azerbaiyan"
$ws.Cells.Item(58,2).Value = 25113
$ws.Cells.Item(58,4).Value = 16150
$ws.Cells.Item(58,5).Value = 8644
$ws.Cells.Item(58,8).Value = 319

# Row 64: 'Nepal' -> 'Nepal'
$ws.Cells.Item(64,2).Value = 17177
$ws.Cells.Item(64,3).Value = 116
$ws.Cells.Item(64,4).Value = 11025
$ws.Cells.Item(64,5).Value = 6113
$ws.Cells.Item(64,7).Value = 1
$ws.Cells.Item(64,8).Value = 39

# Row 65: 'Marruecos' -> 'Marruecos'
$ws.Cells.Item(65,2).Value = 16181
$ws.Cells.Item(65,3).Value = 84
$ws.Cells.Item(65,4).Value = 13666
$ws.Cells.Item(65,5).Value = 2257
$ws.Cells.Item(65,7).Value = 1
$ws.Cells.Item(65,8).Value = 258

# Row 81: 'Republica de Macedonia' -> 'Senegal'
$ws.Cells.Item(81,1).Value = "Senegal"
$ws.Cells.Item(81,2).Value = 8369
$ws.Cells.Item(81,3).Value = 126
$ws.Cells.Item(81,4).Value = 5605
$ws.Cells.Item(81,5).Value = 2611
$ws.Cells.Item(81,7).Value = 3
$ws.Cells.Item(81,8).Value = 153

# Row 82: 'Senegal' -> 'Republica de Macedonia'
$ws.Cells.Item(82,1).Value = "Republica de Macedonia"
$ws.Cells.Item(82,2).Value = 8332
$ws.Cells.Item(82,4).Value = 4468
$ws.Cells.Item(82,5).Value = 3475
$ws.Cells.Item(82,8).Value = 389

# Row 86: 'Finlandia' -> 'Finlandia'
$ws.Cells.Item(86,2).Value = 7296
$ws.Cells.Item(86,4).Value = 6880
$ws.Cells.Item(86,5).Value = 88
$ws.Cells.Item(86,8).Value = 328

# Row 94: 'Mauritania' -> 'Madagascar'
$ws.Cells.Item(94,1).Value = "Madagascar"
$ws.Cells.Item(94,2).Value = 5605
$ws.Cells.Item(94,3).Value = 262
$ws.Cells.Item(94,4).Value = 2811
$ws.Cells.Item(94,5).Value = 2751
$ws.Cells.Item(94,7).Value = 4
$ws.Cells.Item(94,8).Value = 43

# Row 95: 'Madagascar' -> 'Mauritania'
$ws.Cells.Item(95,1).Value = "Mauritania"
$ws.Cells.Item(95,2).Value = 5518
$ws.Cells.Item(95,4).Value = 2664
$ws.Cells.Item(95,5).Value = 2707
$ws.Cells.Item(95,8).Value = 147

# Row 102: 'Albania' -> 'Albania'
$ws.Cells.Item(102,2).Value = 3752
$ws.Cells.Item(102,3).Value = 85
$ws.Cells.Item(102,4).Value = 2091
$ws.Cells.Item(102,5).Value = 1560
$ws.Cells.Item(102,7).Value = 4
$ws.Cells.Item(102,8).Value = 101

# Row 140: 'Uganda' -> 'Uganda'
$ws.Cells.Item(140,2).Value = 1043
$ws.Cells.Item(140,3).Value = 3
$ws.Cells.Item(140,4).Value = 1004
$ws.Cells.Item(140,5).Value = 39

# Row 155: 'Malta' -> 'Malta'
$ws.Cells.Item(155,4).Value = 661
$ws.Cells.Item(155,5).Value = 4

# Row 176: 'Gibraltar' -> 'Gibraltar'
$ws.Cells.Item(176,4).Value = 180
$ws.Cells.Item(176,5).Value = 0

# Row 209: 'Islas Malvinas' -> 'Groenlandia'
$ws.Cells.Item(209,1).Value = "Groenlandia"

# Row 210: 'Groenlandia' -> 'Islas Malvinas'
$ws.Cells.Item(210,1).Value = "Islas Malvinas"

# Row 215: 'Islas Virgenes Britanicas' -> 'Bonaire, San Eustaquio y Saba'
$ws.Cells.Item(215,1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(215,2).Value = 9
$ws.Cells.Item(215,3).Value = 2
$ws.Cells.Item(215,5).Value = 2
$ws.Cells.Item(215,8).Value = 0

# Row 216: 'Bonaire, San Eustaquio y Saba' -> 'Islas Virgenes Britanicas'
$ws.Cells.Item(216,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(216,2).Value = 8
$ws.Cells.Item(216,8).Value = 1
